# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "537.14", "1.00").
# Force each cell to text format first so Excel keeps the literal string
# (matching the existing inline-string cells) instead of silently parsing
# it as a number and dropping significant trailing/leading zeros. (Setting
# NumberFormat on a multi-area Range does not reliably apply per-cell, so
# this is done one cell at a time.)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.402.55"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.154.84"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.14"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.63"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.33"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.699.68"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.74"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000169"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.484.24"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.144.51"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.25"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.85"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.78"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.11"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.514"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.14"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0867"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.92"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.67"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.22"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.30"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.646.84"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0680"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.20"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.707"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.200.14"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.21"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.978"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.24"

# Drop the temporary text-number-format override so the cells end up with
# no explicit style, same as the untouched price cells around them.
$ws.Range("D2").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()

# Column E holds percentage-change text (e.g. "  +2.25%  "); Excel leaves
# these alone as text because of the leading/trailing spaces, so no special
# handling is required.
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +9.42%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +3.76%  "
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("E15").Value = "  +6.66%  "
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("E18").Value = "  +6.57%  "
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("E20").Value = "  +6.42%  "
$ws.Range("E21").Value = "  +8.39%  "
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +14.56%  "
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("E31").Value = "  +6.32%  "
$ws.Range("E32").Value = "  +5.19%  "
$ws.Range("E33").Value = "  +8.28%  "
$ws.Range("E34").Value = "  +5.63%  "
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  +4.52%  "
$ws.Range("E37").Value = "  +12.81%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  +10.11%  "
$ws.Range("E40").Value = "  +6.23%  "
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("E42").Value = "  +4.78%  "
$ws.Range("E43").Value = "  +6.02%  "
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("E45").Value = "  +8.54%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("E48").Value = "  +11.20%  "
$ws.Range("E49").Value = "  +4.59%  "
$ws.Range("E50").Value = "  +5.36%  "
$ws.Range("E51").Value = "  +5.28%  "

